# feat: add 2022-Q1 data
#
# The existing "总计" (totals) sheet is renamed to "2022-Q1" and re-populated
# with the per-fund holding detail for the new quarter (mirroring the layout
# already used by the "2021-Q3" / "2021-Q4" sheets). A brand new "总计" sheet
# is then inserted right after it, carrying forward the previous totals plus
# a new row summarising 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Turn the old "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Wipe the old totals-table content/format so we can rebuild the wider
# fund-detail table from scratch.
$q1.Range("A1:H10").Clear()

# Borrow the header / row-index styling (bold font + border) from the
# sibling quarter sheet so the new table matches the rest of the workbook.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Rows = @(
    @(0, "006002", "工银瑞信医药健康股票A",        "34.28", "87.97", "1.80", "0.6170", 10),
    @(1, "007110", "国投瑞银港股通价值发现混合",     "23.33", "93.33", "2.64", "0.6159", 8),
    @(2, "010088", "工银瑞信优质成长混合A",          "19.41", "82.44", "2.27", "0.4406", 8),
    @(3, "006003", "工银瑞信医药健康股票C",          "10.03", "87.97", "1.80", "0.1805", 10),
    @(4, "470888", "汇添富香港优势精选混合 (QDII)",  "1.38",  "67.74", "3.00", "0.0414", 7),
    @(5, "010089", "工银瑞信优质成长混合C",          "1.34",  "82.44", "2.27", "0.0304", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Range("A${r}").Value = $row[0]
    $q1.Range("B${r}").NumberFormat = "@"
    $q1.Range("B${r}").Value = $row[1]
    $q1.Range("B${r}").Style = "Normal"
    $q1.Range("C${r}").Value = $row[2]
    $dataRange = $q1.Range("D${r}:G${r}")
    $dataRange.NumberFormat = "@"
    $q1.Range("D${r}").Value = $row[3]
    $q1.Range("E${r}").Value = $row[4]
    $q1.Range("F${r}").Value = $row[5]
    $q1.Range("G${r}").Value = $row[6]
    $dataRange.Style = "Normal"
    $q1.Range("H${r}").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Insert a fresh "总计" sheet right after "2022-Q1" with the rolled-up
#    per-quarter summary (previous rows + the new 2022-Q1 row on top).
#    Cloning the template sheet (instead of Worksheets.Add) keeps the same
#    sheet-level setup (outline/page properties) used by every other tab.
# ---------------------------------------------------------------------------
$template.Copy($null, $q1)
$total = $wb.Worksheets.Item(4)
$total.Name = "总计"
$total.Range("A1:H10").Clear()

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 1.93

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 1.43

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.74

# Restore the originally active tab (sheet 1) now that all edits are done.
$wb.Worksheets.Item(1).Activate()
